$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 396, shifting existing rows 396..483 down to 397..484
$newRow = $ws.Rows.Item(396)
$newRow.Insert(-4121)  # xlShiftDown

# Populate the newly inserted row 396 with the new record's data
$ws.Cells.Item(396, 1).Value = 4
$ws.Cells.Item(396, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(396, 3).Value = "Los Lagos"
$ws.Cells.Item(396, 4).Value = 44964
$ws.Cells.Item(396, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(396, 5).Value = 10
$ws.Cells.Item(396, 6).Value = 100112008
$ws.Cells.Item(396, 7).Value = "Coliflor"
$ws.Cells.Item(396, 8).Value = "Sin especificar"
$ws.Cells.Item(396, 9).Value = "Primera"
$ws.Cells.Item(396, 10).Value = 1200
$ws.Cells.Item(396, 11).Value = 1800
$ws.Cells.Item(396, 12).Value = 1800
$ws.Cells.Item(396, 13).Value = 1800
$ws.Cells.Item(396, 14).Value = "`$/unidad"
$ws.Cells.Item(396, 15).Value = "Región Metropolitana"
$ws.Cells.Item(396, 16).Value = 1800
$ws.Cells.Item(396, 17).Value = 1
$ws.Cells.Item(396, 18).Value = "Hortaliza"
